$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Unveiling the Mysteries of the Quantum World" "A Journey Through the Magic of Mathematics"

# --- Author name line (collapses the "Dr" / "." / " Isabella Santos" runs) ---
Replace-Text "Dr. Isabella Santos" "Samuel Davies"

# --- E-mail line ---
# "isabella" -> "samueldavies@edumail" (own run)
Replace-Text "isabella" "samueldavies@edumail"
# "santos@eliteuniversity" -> "org" ; drop the trailing ".edu" but keep the
# "." run that separates the two halves of the address.
Replace-Text "santos@eliteuniversity" "org"
Replace-Text "org.edu" "org"

# --- Body paragraph, first blank-line-delimited block ---
$oldA = "In the intricate realm of quantum physics, where the fabric of reality assumes an enigmatic and ethereal character, we embark on an exploration of the fundamental nature of matter and energy. This fascinating journey leads us to the heart of quantum mechanics, where the laws governing the behavior of particles at the atomic and subatomic levels diverge markedly from those governing the macroscopic world we perceive with our senses. In this realm, particles can exist in multiple states simultaneously, phenomena such as entanglement defy classical notions of locality, and the very act of observation influences the outcome of experiments."
$newA = "Mathematics, the universal language of the universe, beckons us on an awe-inspiring intellectual adventure. This enchanting realm has captured the imagination of brilliant minds throughout history, inspiring discoveries that have shaped our understanding of the world. From ancient civilizations to modern-day frontiers, mathematics weaves intricate threads connecting science, technology, engineering, and art. It is the orchestra conductor of our universe, orchestrating the rhythm and melody of existence."
Replace-Text $oldA $newA

# --- Body paragraph, second blank-line-delimited block ---
$oldB = "As we delve deeper into the quantum realm, we confront paradoxes that challenge our understanding of reality. Schrodinger's cat, a thought experiment conceived by Austrian physicist Erwin Schrodinger, epitomizes this paradoxical nature, presenting a scenario where a cat's fate is entangled with the outcome of a random quantum event. This experiment highlights the perplexing superposition principle, where particles can occupy multiple states until they are observed, challenging our intuitive understanding of the world."
$newB = "In the symphony of mathematics, we witness the harmony of patterns, the elegance of symmetry, and the power of logical reasoning. It holds the key to unlocking nature's deepest secrets, empowering us to unveil the enigmas of the cosmos. Mathematics empowers us to navigate the complexity of modern life, from financial transactions to intricate engineering marvels. It is the language of innovation, unraveling new frontiers of knowledge and shaping the trajectory of human progress."
Replace-Text $oldB $newB

# --- Body paragraph, third blank-line-delimited block ---
$oldC = "The study of quantum mechanics has not only illuminated the intricate behavior of particles at the foundation of matter but also yielded practical applications with far-reaching implications. Quantum technologies, such as quantum computing and quantum cryptography, hold the potential to revolutionize industries, offering unprecedented computational power and unbreakable communication channels. Quantum mechanics forms the bedrock of modern physics, profoundly influencing fields as diverse as cosmology, where it paves the way for understanding the origin and evolution of the universe, and condensed matter physics, where it guides the design of novel materials with extraordinary properties."
$newC = "Mathematics is more than a subject; it's a way of thinking, a lens through which we can perceive the world with greater clarity. It cultivates analytical and problem-solving skills, nurturing critical thinking and inspiring creativity. By delving into the depths of mathematics, we cultivate a mindset that embraces precision, logic, and imagination, enabling us to become more effective problem solvers and informed decision-makers."
Replace-Text $oldC $newC

# --- Summary paragraph ---
$oldS = "Our exploration of the quantum realm has unveiled a universe governed by laws that transcend our classical intuitions and challenge our understanding of reality. Quantum mechanics, with its superposition principle, entanglement, and paradoxical nature, has revolutionized our comprehension of the fundamental constituents of matter and energy. This profound insight has not only expanded our knowledge of the universe but has also given rise to cutting-edge technologies with the potential to transform industries and reshape society. The study of quantum mechanics continues to push the boundaries of human understanding, promising further breakthroughs and a deeper appreciation of the cosmos."
$newS = "This exploration of mathematics highlights its remarkable allure, unveiling its role as the universal language underpinning the fabric of our universe. Mathematics weaves enchanting connections between the worlds of science, technology, and art, inspiring discoveries that have profoundly shaped our understanding of existence. It empowers us to unlock nature's secrets, navigate the complexities of modern life, and cultivate analytical minds capable of solving real-world problems. Embracing the enchantment of mathematics, we transform into informed thinkers and effective problem-solvers, ready to contribute to a future shaped by innovation and discovery."
Replace-Text $oldS $newS

# --- Trailing empty paragraph added at the end of the document body ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
